$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Epoch Accuracy (column B) values per the recorded training run
$ws.Range("B3").Value = 0.375
$ws.Range("B4").Value = 0.3125
$ws.Range("B7").Value = 0.328125
$ws.Range("B9").Value = 0.3125
$ws.Range("B10").Value = 0.25
$ws.Range("B11").Value = 0.28125
$ws.Range("B12").Value = 0.296875
$ws.Range("B13").Value = 0.25
$ws.Range("B14:B15").Value = 0.28125
$ws.Range("B16:B17").Value = 0.265625
$ws.Range("B18").Value = 0.28125
$ws.Range("B19").Value = 0.296875
$ws.Range("B20").Value = 0.265625
$ws.Range("B22").Value = 0.25
$ws.Range("B23").Value = 0.265625
$ws.Range("B24:B27").Value = 0.28125
$ws.Range("B30").Value = 0.28125
$ws.Range("B32").Value = 0.265625
$ws.Range("B33").Value = 0.359375
$ws.Range("B34").Value = 0.265625
$ws.Range("B35").Value = 0.34375
$ws.Range("B36").Value = 0.296875
$ws.Range("B37").Value = 0.265625
$ws.Range("B38").Value = 0.25
$ws.Range("B39:B40").Value = 0.265625
$ws.Range("B41").Value = 0.25
$ws.Range("B42:B43").Value = 0.234375
$ws.Range("B44").Value = 0.21875
$ws.Range("B45:B102").Value = 0.203125
$ws.Range("B103").Value = 0.21875
$ws.Range("B104").Value = 0.125
$ws.Range("B105").Value = 0.203125
$ws.Range("B106").Value = 0.25
$ws.Range("B107").Value = 0.140625
$ws.Range("B109").Value = 0.109375
$ws.Range("B110").Value = 0.078125
$ws.Range("B111").Value = 0.25
$ws.Range("B112").Value = 0.09375
$ws.Range("B114").Value = 0.078125
$ws.Range("B115").Value = 0.109375
$ws.Range("B116:B117").Value = 0.15625
$ws.Range("B118").Value = 0.180327868852459

# Refresh the DisplayOutputs repr text (new object instance/address from the latest run)
$ws.Range("A102:A118").Value = "<__main__.DisplayOutputs object at 0x7f77e8408610>"
